# Insert two new "Author" paragraphs ("Ben Jarman" and "Helen Fair")
# immediately after the Subtitle paragraph ("Laws, policies, and
# practical realities") and before the Date paragraph ("2024-07-10").

$d = $word.ActiveDocument

# Find the Subtitle paragraph.
$subtitlePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Subtitle") {
        $subtitlePara = $p
        break
    }
}

# The paragraph right after the subtitle (currently the Date paragraph)
# is where the new content needs to be inserted in front of.
$nextPara = $subtitlePara.Next()
$insertRange = $d.Range($nextPara.Range.Start, $nextPara.Range.Start)

# Build a WordprocessingML package fragment containing the two new
# "Author" paragraphs. A trailing empty "Author" paragraph is included
# because InsertXML merges the final paragraph mark of the inserted
# fragment into the paragraph at the insertion point (the same way
# pasting multiple paragraphs works in Word); that trailing paragraph
# is removed again below once it has served its purpose of keeping
# "Helen Fair" on its own paragraph.
$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Author"/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">Ben Jarman</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Author"/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">Helen Fair</w:t>
            </w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Author"/>
            </w:pPr>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$null = $insertRange.InsertXML($xml)

# Re-resolve the Date paragraph: the old $nextPara reference is stale
# after InsertXML has changed the document, since it still refers to
# its original (now stale) text range rather than tracking the moved
# paragraph. Remove the spurious empty "Author" paragraph left
# immediately before the (re-resolved) Date paragraph.
$datePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Date") {
        $datePara = $p
        break
    }
}
$null = $datePara.Previous().Range.Delete()
